$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data: mesas de examen added for "Enfermería Universitaria" (rows 76-82)
# and "Licenciatura En Enfermería" (rows 83-85).
# ---------------------------------------------------------------------------

$rows = @(
    @{ r=76; B="Enfermería Universitaria"; C="Herramientas De Informatica";    D="Zuñiga";      E="Pezzini";       F="Constanzo M."; G="Farias R."; H=43367; I=0.41666666666666669 },
    @{ r=77; B="Enfermería Universitaria"; C="Enfermeria Basica";              D="Rojas A.";     E="Mansilla V.";   F="Buera S.";      G="Brito E.";  H=43371; I=0.70833333333333337 },
    @{ r=78; B="Enfermería Universitaria"; C="Psicologia Evolutiva";           D="Mansilla S.";  E="Rosales K.";    F="Navarro F.";    G="Rossi V.";  H=43369; I=0.375 },
    @{ r=79; B="Enfermería Universitaria"; C="Ciencia Universidad Y Sociedad"; D="Jaremchuk";    E="Vilaboa";       F="Carcamo Y.";    G="Musci C."; H=43370; I=0.625 },
    @{ r=80; B="Enfermería Universitaria"; C="Antropologia Sociocultural";     D="Alvarez P.";   E="Campan";        F="Enrici";        G=$null;      H=43368; I=0.625 },
    @{ r=81; B="Enfermería Universitaria"; C="Ciencias Biologicas";            D="Navarro O.";   E="Brandoni";      F="Miro";          G=$null;      H=43371; I=0.75 },
    @{ r=82; B="Enfermería Universitaria"; C="Analisis Y Producción Del Discurso"; D="Dalla Costa"; E="Guidetti C."; F="Altamirano"; G="Bahamonde S."; H=43369; I=0.625 },
    @{ r=83; B="Licenciatura En Enfermería"; C="Enfermeria En Alto Riesgo";    D="Vega P.";      E="Soria D.";      F="Galarza M.";    G=$null;      H=$null; I=$null },
    @{ r=84; B="Licenciatura En Enfermería"; C="Epidemiologia Y Estadistica";  D="Firnkorn M.";  E="Vallejos P.";   F=$null;           G=$null;      H=$null; I=$null },
    @{ r=85; B="Licenciatura En Enfermería"; C="Investigacion En Enfermeria";  D="Ojeda S.";      E="Mansilla V.";   F="Oyarzo V.";     G=$null;      H=$null; I=$null }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    if ($row.D -ne $null) { $ws.Cells.Item($r, 4).Value = $row.D }
    if ($row.E -ne $null) { $ws.Cells.Item($r, 5).Value = $row.E }
    if ($row.F -ne $null) { $ws.Cells.Item($r, 6).Value = $row.F }
    if ($row.G -ne $null) { $ws.Cells.Item($r, 7).Value = $row.G }
    if ($row.H -ne $null) { $ws.Cells.Item($r, 8).Value = $row.H }
    if ($row.I -ne $null) { $ws.Cells.Item($r, 9).Value = $row.I }

    # Carry over the same number/visual formatting the rest of the table uses
    # (column D:G -> "s=3" font style, H -> short date, I -> time) by copying
    # formats from row 71, which already uses that exact style combination.
    $ws.Range("D71").Copy()
    $ws.Range("D" + $r).PasteSpecial(-4122)
    $ws.Range("E71").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)
    $ws.Range("F71").Copy()
    $ws.Range("F" + $r).PasteSpecial(-4122)
    if ($row.G -ne $null) {
        $ws.Range("G71").Copy()
        $ws.Range("G" + $r).PasteSpecial(-4122)
    }
    if ($row.H -ne $null) {
        $ws.Range("H71").Copy()
        $ws.Range("H" + $r).PasteSpecial(-4122)
    }
    if ($row.I -ne $null) {
        $ws.Range("I71").Copy()
        $ws.Range("I" + $r).PasteSpecial(-4122)
    }
}

# Row 82 & onward's B column loses the extra font style (matches source diff:
# B82/B83/B84/B85 have no explicit style while B76:B81 do).
$ws.Range("B71").Copy()
$ws.Range("B76").PasteSpecial(-4122)
$ws.Range("B77").PasteSpecial(-4122)
$ws.Range("B78").PasteSpecial(-4122)
$ws.Range("B79").PasteSpecial(-4122)
$ws.Range("B80").PasteSpecial(-4122)
$ws.Range("B81").PasteSpecial(-4122)

# Sheet view / selection moved as part of this edit.
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("B85").Select()
